$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 219.72
$ws.Range("I15").Value = 219.72
$ws.Range("K15").Value = 659.16
$ws.Range("M15").Value = -490.16
$ws.Range("H138").Value = 3516.6667
$ws.Range("I138").Value = 1704.6451
$ws.Range("J138").Value = 4793.3184
$ws.Range("K138").Value = 5113.9353
$ws.Range("L138").Value = 14379.9552
$ws.Range("M138").Value = 26.0646999999999
$ws.Range("N138").Value = -24659.9552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14880.93
$ws.Range("I32").Value = 11298.963
$ws.Range("K32").Value = 11298.963
$ws.Range("M32").Value = -11011.963
$ws.Range("H122").Value = 2568980.5
$ws.Range("I122").Value = 2853978.2
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8561934.600000001
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8559484.600000001
$ws.Range("N122").Value = -16900
$ws.Range("H137").Value = 36157.5
$ws.Range("J137").Value = 36157.5
$ws.Range("L137").Value = 36157.5
$ws.Range("N137").Value = -46357.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 300.36365
$ws.Range("J64").Value = 288
$ws.Range("L64").Value = 288
$ws.Range("N64").Value = -738
$ws.Range("H67").Value = 300.36365
$ws.Range("J67").Value = 288
$ws.Range("L67").Value = 288
$ws.Range("N67").Value = -1848
$ws.Range("H94").Value = 1146.7368
$ws.Range("I94").Value = 969.1429000000001
$ws.Range("J94").Value = 1644
$ws.Range("K94").Value = 969.1429000000001
$ws.Range("L94").Value = 1644
$ws.Range("M94").Value = -518.1429000000001
$ws.Range("N94").Value = -2546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15159672
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15159672
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15159672
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -15160262
$ws.Range("H34").Value = 15159672
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 15159672
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 15159672
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -15160076
$ws.Range("H58").Value = 6544940
$ws.Range("I58").Value = 7577448.5
$ws.Range("J58").Value = 2001903
$ws.Range("K58").Value = 7577448.5
$ws.Range("L58").Value = 2001903
$ws.Range("M58").Value = -7577245.5
$ws.Range("N58").Value = -2002309
$ws.Range("H122").Value = 4291.3184
$ws.Range("I122").Value = 3710.2307
$ws.Range("J122").Value = 5130.6665
$ws.Range("K122").Value = 11130.6921
$ws.Range("L122").Value = 15391.9995
$ws.Range("M122").Value = -8680.6921
$ws.Range("N122").Value = -20291.9995
$ws.Range("H136").Value = 6544940
$ws.Range("I136").Value = 7577448.5
$ws.Range("J136").Value = 2001903
$ws.Range("K136").Value = 22732345.5
$ws.Range("L136").Value = 6005709
$ws.Range("M136").Value = -22729795.5
$ws.Range("N136").Value = -6010809

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2500
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 6000
$ws.Range("M80").Value = -5064
$ws.Range("H83").Value = 2500
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 18000
$ws.Range("M83").Value = -13320
$ws.Range("H106").Value = 3887.5
$ws.Range("J106").Value = 3887.5
$ws.Range("L106").Value = 11662.5
$ws.Range("N106").Value = -13554.5
$ws.Range("H109").Value = 1092.875
$ws.Range("I109").Value = 963.2857
$ws.Range("J109").Value = 2000
$ws.Range("K109").Value = 2889.8571
$ws.Range("L109").Value = 6000
$ws.Range("M109").Value = -1849.8571
$ws.Range("N109").Value = -8080
$ws.Range("H110").Value = 6853.3125
$ws.Range("I110").Value = 4513.25
$ws.Range("J110").Value = 7633.3335
$ws.Range("K110").Value = 13539.75
$ws.Range("L110").Value = 22900.0005
$ws.Range("M110").Value = -9449.75
$ws.Range("N110").Value = -31080.0005
$ws.Range("H112").Value = 3010.9375
$ws.Range("I112").Value = 1350
$ws.Range("J112").Value = 3476
$ws.Range("K112").Value = 4050
$ws.Range("L112").Value = 10428
$ws.Range("M112").Value = -2942
$ws.Range("N112").Value = -12644
$ws.Range("H118").Value = 1483
$ws.Range("J118").Value = 2963.3333
$ws.Range("L118").Value = 8889.999899999999
$ws.Range("N118").Value = -11375.9999
$ws.Range("H121").Value = 722.1818
$ws.Range("I121").Value = 358.75
$ws.Range("J121").Value = 929.8570999999999
$ws.Range("K121").Value = 1076.25
$ws.Range("L121").Value = 2789.5713
$ws.Range("M121").Value = 233.75
$ws.Range("N121").Value = -5409.5713
$ws.Range("H123").Value = 6698.75
$ws.Range("I123").Value = 3647.5
$ws.Range("J123").Value = 9750
$ws.Range("K123").Value = 10942.5
$ws.Range("L123").Value = 29250
$ws.Range("M123").Value = -8492.5
$ws.Range("N123").Value = -34150
$ws.Range("H129").Value = 1235.375
$ws.Range("I129").Value = 1023.0769
$ws.Range("J129").Value = 1486.2727
$ws.Range("K129").Value = 3069.2307
$ws.Range("L129").Value = 4458.8181
$ws.Range("M129").Value = 1930.7693
$ws.Range("N129").Value = -14458.8181
$ws.Range("H131").Value = 2041861.4
$ws.Range("I131").Value = 6667363.5
$ws.Range("J131").Value = 1198.7646
$ws.Range("K131").Value = 20002090.5
$ws.Range("L131").Value = 3596.2938
$ws.Range("M131").Value = -19997050.5
$ws.Range("N131").Value = -13676.2938
$ws.Range("H133").Value = 63288.832
$ws.Range("I133").Value = 131943.75
$ws.Range("J133").Value = 8364.9
$ws.Range("K133").Value = 395831.25
$ws.Range("L133").Value = 25094.7
$ws.Range("M133").Value = -390771.25
$ws.Range("N133").Value = -35214.7
$ws.Range("H137").Value = 17117.04
$ws.Range("I137").Value = 10466.154
$ws.Range("J137").Value = 24322.166
$ws.Range("K137").Value = 31398.462
$ws.Range("L137").Value = 72966.49800000001
$ws.Range("M137").Value = -26298.462
$ws.Range("N137").Value = -83166.49800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 6779.4707
$ws.Range("J57").Value = 6779.4707
$ws.Range("L57").Value = 6779.4707
$ws.Range("N57").Value = -8419.4707
$ws.Range("H126").Value = 12056.158
$ws.Range("I126").Value = 15226.571
$ws.Range("J126").Value = 3179
$ws.Range("K126").Value = 45679.713
$ws.Range("L126").Value = 9537
$ws.Range("M126").Value = -43209.713
$ws.Range("N126").Value = -14477
$ws.Range("H132").Value = 10419358
$ws.Range("I132").Value = 15153703
$ws.Range("J132").Value = 3798.6
$ws.Range("K132").Value = 45461109
$ws.Range("L132").Value = 11395.8
$ws.Range("M132").Value = -45458579
$ws.Range("N132").Value = -16455.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3260493.8
$ws.Range("I122").Value = 3971736
$ws.Range("J122").Value = 1431585
$ws.Range("K122").Value = 11915208
$ws.Range("L122").Value = 4294755
$ws.Range("M122").Value = -11912758
$ws.Range("N122").Value = -4299655

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 66667510
$ws.Range("I107").Value = 100000770
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 300002310
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -300000390
$ws.Range("N107").Value = -6840
$ws.Range("H136").Value = 2978285.2
$ws.Range("I136").Value = 1883.3334
$ws.Range("J136").Value = 11907491
$ws.Range("K136").Value = 5650.0002
$ws.Range("L136").Value = 35722473
$ws.Range("M136").Value = -3100.0002
$ws.Range("N136").Value = -35727573
